$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) store these numeric-looking values as
# plain text in the workbook. A leading apostrophe forces Excel to keep them
# as text instead of auto-converting to a number/percentage.

$ws.Range("D2").Value = "'303.27"
$ws.Range("E2").Value = "'-4.90%"
$ws.Range("D3").Value = "'35.17"
$ws.Range("E3").Value = "'-2.67%"
$ws.Range("D4").Value = "'5.063"
$ws.Range("E4").Value = "'-2.83%"
$ws.Range("D5").Value = "'0.08000"
$ws.Range("E5").Value = "'-2.82%"
$ws.Range("D6").Value = "'1.931"
$ws.Range("E6").Value = "'-10.32%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'7.749"
$ws.Range("E7").Value = "'-3.93%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.905"
$ws.Range("E8").Value = "'3.70%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9207"
$ws.Range("E9").Value = "'-0.69%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1228"
$ws.Range("E10").Value = "'20.21%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1850"
$ws.Range("E11").Value = "'-2.04%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09451"
$ws.Range("E12").Value = "'3.20%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03588"
$ws.Range("E13").Value = "'-0.98%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09857"
$ws.Range("E14").Value = "'-0.67%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001389"
$ws.Range("E15").Value = "'-3.39%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005750"
$ws.Range("E16").Value = "'1.67%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.499"
$ws.Range("E17").Value = "'1.08%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.057"
$ws.Range("E18").Value = "'-1.88%"
$ws.Range("E19").Value = "'2.11%"
$ws.Range("D20").Value = "'0.1282"
$ws.Range("E20").Value = "'-1.42%"
$ws.Range("D21").Value = "'5.030"
$ws.Range("E21").Value = "'-0.62%"
$ws.Range("E22").Value = "'12.56%"
$ws.Range("D23").Value = "'0.04508"
$ws.Range("E23").Value = "'-1.91%"
$ws.Range("E24").Value = "'-2.39%"
$ws.Range("D25").Value = "'0.004852"
$ws.Range("E25").Value = "'2.52%"
$ws.Range("E26").Value = "'-0.11%"
$ws.Range("E27").Value = "'-6.87%"
$ws.Range("D39").Value = "'0.01933"
$ws.Range("E39").Value = "'-3.55%"
$ws.Range("D40").Value = "'0.04743"
$ws.Range("E40").Value = "'-4.37%"
$ws.Range("D41").Value = "'0.007559"
$ws.Range("E41").Value = "'-3.23%"
$ws.Range("D42").Value = "'0.009538"
$ws.Range("E42").Value = "'25.72%"
$ws.Range("D43").Value = "'0.1330"
$ws.Range("E43").Value = "'-4.98%"
$ws.Range("E44").Value = "'0.57%"
$ws.Range("D45").Value = "'0.01114"
$ws.Range("E45").Value = "'-6.09%"
$ws.Range("D46").Value = "'0.00006293"
$ws.Range("E46").Value = "'-2.43%"
$ws.Range("E47").Value = "'-0.10%"
$ws.Range("E48").Value = "'57.11%"
$ws.Range("E49").Value = "'-31.39%"
$ws.Range("E50").Value = "'-0.10%"
$ws.Range("E51").Value = "'-0.10%"
